# Insert a new weekly price record as the new row 7 of the "Orégano" sheet.
# This pushes the previously-existing rows 7..123 down to 8..124, which
# matches the author's commit ("Fruta / hortaliza, semanal") that prepends
# the latest week's observation to the top of the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7 and below down by one row.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 45245
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = 100112029
$ws.Range("G7").Value = "Orégano"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 16
$ws.Range("K7").Value = 21000
$ws.Range("L7").Value = 21000
$ws.Range("M7").Value = 21000
$ws.Range("N7").Value = "$/docena de atados"
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 7000
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = "Hortaliza"
